$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Delete row 165 ("「図形の中に三角形はいくつ？」" quiz entry) - all subsequent
# rows shift up by one automatically.
$ws.Rows.Item(165).Delete()
